$d = $word.ActiveDocument

# Color used throughout this revision for the section headings and the
# document header title/date field: RGB(0x02,0x44,0x42) -> OLE BGR long.
$themeColor = 4342786   # 0x024442 in RGB(r,g,b) OLE encoding

# ---------------------------------------------------------------------
# Body headings: add the new font color to the Heading1 / Heading2
# paragraphs ("Calculation Details" and the four "Step N: ..." titles).
# Setting Font.Color on the paragraph Range colors both the paragraph
# mark (pPr/rPr) and the run(s) of text, matching the target markup.
# ---------------------------------------------------------------------
$headingTexts = @(
    "Calculation Details",
    "Step 1: Household Composition",
    "Step 2: Gross Income Test",
    "Step 3: Net Income",
    "Step 4: Final Determination"
)

$total = $d.Paragraphs.Count
for ($i = 1; $i -le $total; $i++) {
    $para = $d.Paragraphs($i)
    $text = $para.Range.Text.Trim()
    if ($headingTexts -contains $text) {
        $para.Range.Font.Color = $themeColor
    }
}

# ---------------------------------------------------------------------
# Page header: the title/date line also gets the new font color,
# covering both the paragraph mark and every run it can reach.
# ---------------------------------------------------------------------
$hdr = $d.Sections(1).Headers(1)
$hdr.Range.Font.Color = $themeColor
